$d = $word.ActiveDocument

# Helper: the XML namespace declaration used for the ad-hoc <w:p> fragments
# we inject via Range.InsertXML to produce truly-empty paragraphs (a
# paragraph that carries only paragraph-mark run properties, with no
# <w:r> child at all - this is what Word leaves behind when you press
# Enter on a blank line without ever typing a character there).
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Add-EmptyParagraph {
    $last = $d.Paragraphs.Last
    $null = $last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $xml = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
    $null = $newPara.Range.InsertXML($xml)
}

function Add-TextParagraph {
    param([string]$text)
    $last = $d.Paragraphs.Last
    $null = $last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $null = $newPara.Range.InsertAfter($text)
}

# Three blank lines after "B6: chạy lệnh npm start như B3"
Add-EmptyParagraph
Add-EmptyParagraph
Add-EmptyParagraph

Add-TextParagraph "Tài khoản Admin : tinpham1510"
Add-TextParagraph "Mật khẩu: 12345678"

Add-EmptyParagraph

Add-TextParagraph "Tài khoản khách hàng : hiep"
Add-TextParagraph "Mật khẩu: 12345678"

# Final paragraph holds a single literal space character
Add-TextParagraph " "
